# Append new "building work study" log entries to the bottom of the
# "Sheet" worksheet (rows 113-159), continuing the existing chat-log
# style rows already present in A1:A112.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$values = @(
    '19 10:56>>> CFD893A460   Gary Tsai',
    '19 16:19>>> CFD893A460   Gary Tsai',
    '19 16:20>>> CFD893A460   Gary Tsai',
    '19 16:43>>> CFD893A460   Gary Tsai',
    '19 16:47>>> CFD893A460   Gary Tsai',
    '19 16:50>>> CFD893A460   Gary Tsai',
    '19 16:52>>> CFD893A460   Gary Tsai',
    '19 17:00>>> CFD893A460   Gary Tsai',
    '19 17:02>>> CFD893A460   Gary Tsai',
    '19 17:17>>> CFD893A460   Gary Tsai',
    '19 17:20>>> CFD893A460   Gary Tsai',
    '19 17:20>>> CFD893A460   Gary Tsai',
    '19 17:26>>> CFD893A460   Gary Tsai',
    '19 17:28>>> CFD893A460   Gary Tsai',
    '19 17:35>>> CFD893A460   Gary Tsai',
    '19 17:44>>> CFD893A460   Gary Tsai',
    '19 17:44>>> CFD893A460   Gary Tsai',
    '19 17:44>>> 0FD8AD9C60   everton mendozaholmes',
    '19 17:46>>> 0FD8AD9C60   everton mendozaholmes',
    '19 17:46>>> 0FD8AD9C60   everton mendozaholmes',
    '19 17:46>>> CFD893A460   Gary Tsai',
    '19 18:00>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:00>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:00>>> CFD893A460   Gary Tsai',
    '19 18:09>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:09>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:10>>> CFD893A460   Gary Tsai',
    '19 18:19>>> CFD893A460   Gary Tsai',
    '19 18:22>>> CFD893A460   Gary Tsai',
    '19 18:22>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:22>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:25>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:27>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:35>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:37>>> CFD893A460   Gary Tsai',
    '19 18:37>>> CFD893A460   Gary Tsai',
    '19 18:37>>> CFD893A460   Gary Tsai',
    '19 18:38>>> CFD893A460   Gary Tsai',
    '19 18:38>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:41>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:44>>> 0FD8AD9C60   everton mendozaholmes',
    '19 18:53>>> CFD893A460   Gary Tsai',
    '19 18:54>>> CFD893A460   Gary Tsai',
    '19 18:57>>> CFD893A460   Gary Tsai',
    '19 18:58>>> CFD893A460   Gary Tsai',
    '19 18:59>>> CFD893A460   Gary Tsai',
    '19 19:00>>> CFD893A460   Gary Tsai'
)

$startRow = 113
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

# "Monthly_STAT" sheet content (Month / August 2016 header) is unchanged
# by this edit; only its shared-string indices shift automatically as a
# side effect of the new strings added above.

